$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column D, matching style/formatting of existing headers (A1:C1)
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Tipo"

# Update MSE (B) and R2 (C) values with new computed figures, and fill new D column
$ws.Range("B2").Value = 0.4927202043440829
$ws.Range("C2").Value = 0.9855275168862212
$ws.Range("D2").Value = "multiple"

$ws.Range("B3").Value = 0.09314024093052152
$ws.Range("C3").Value = 0.9987127766512816
$ws.Range("D3").Value = "multiple"

$ws.Range("B4").Value = 0.03285527010160284
$ws.Range("C4").Value = 0.9996613677366777
$ws.Range("D4").Value = "multiple"

$ws.Range("B5").Value = 0.09464358933099511
$ws.Range("C5").Value = 0.9994362351423224
$ws.Range("D5").Value = "multiple"
